$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 511; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

$ws.Rows.Item(511).RowHeight = 15

$ws.Cells.Item(512, 1).Value = "A 45532-2023"
$ws.Cells.Item(512, 2).Value = 45194
$ws.Cells.Item(512, 3).Value = 45202
$ws.Cells.Item(512, 4).Value = "VÄSTRA GÖTALANDS LÄN"
$ws.Cells.Item(512, 5).Value = "ULRICEHAMN"
$ws.Cells.Item(512, 7).Value = 0.5
$ws.Cells.Item(512, 8).Value = 0
$ws.Cells.Item(512, 9).Value = 0
$ws.Cells.Item(512, 10).Value = 0
$ws.Cells.Item(512, 11).Value = 0
$ws.Cells.Item(512, 12).Value = 0
$ws.Cells.Item(512, 13).Value = 0
$ws.Cells.Item(512, 14).Value = 0
$ws.Cells.Item(512, 15).Value = 0
$ws.Cells.Item(512, 16).Value = 0
$ws.Cells.Item(512, 17).Value = 0

# Copy styles from row 511 to row 512 (B,C date style, R wrap style)
$ws.Cells.Item(512, 2).NumberFormat = $ws.Cells.Item(511, 2).NumberFormat
$ws.Cells.Item(512, 3).NumberFormat = $ws.Cells.Item(511, 3).NumberFormat
$ws.Cells.Item(512, 18).WrapText = $ws.Cells.Item(511, 18).WrapText

Write-Host "done"
